$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value that was bumped by one
# day (45188 -> 45189, i.e. 2023-09-19 -> 2023-09-20) for every data row.
$lastRow = $ws.Range("A1").End(4).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
